$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '67.143.17'
$ws.Range('E2').Value = '  -0.98%  '
Set-TextValue $ws.Range('D3') '2.472.36'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('E4').Value = '  +0.05%  '
Set-TextValue $ws.Range('D5') '583.48'
Set-TextValue $ws.Range('D6') '170.22'
$ws.Range('E6').Value = '  -1.73%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -2.16%  '
Set-TextValue $ws.Range('D9') '2.471.92'
$ws.Range('E9').Value = '  -2.75%  '
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('E13').Value = '  -3.77%  '
Set-TextValue $ws.Range('D14') '25.64'
$ws.Range('E14').Value = '  -3.48%  '
Set-TextValue $ws.Range('D15') '2.920.71'
$ws.Range('E15').Value = '  +0.03%  '
Set-TextValue $ws.Range('D16') '66.996.39'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('E17').Value = '  -4.52%  '
Set-TextValue $ws.Range('D18') '2.473.95'
$ws.Range('E18').Value = '  -2.66%  '
Set-TextValue $ws.Range('D19') '11.16'
$ws.Range('E19').Value = '  -5.72%  '
Set-TextValue $ws.Range('D20') '7.61'
$ws.Range('E20').Value = '  -3.54%  '
Set-TextValue $ws.Range('D21') '353.70'
$ws.Range('E21').Value = '  -4.54%  '
Set-TextValue $ws.Range('D22') '4.04'
$ws.Range('E22').Value = '  -2.80%  '
Set-TextValue $ws.Range('D23') '1.00'
$ws.Range('E23').Value = '  +0.02%  '
Set-TextValue $ws.Range('D24') '68.99'
$ws.Range('E24').Value = '  -3.85%  '
Set-TextValue $ws.Range('D25') '4.26'
$ws.Range('E25').Value = '  -7.28%  '
$ws.Range('E26').Value = '  -7.18%  '
Set-TextValue $ws.Range('D27') '9.25'
$ws.Range('E27').Value = '  -7.29%  '
$ws.Range('E28').Value = '  -9.71%  '
Set-TextValue $ws.Range('D29') '2.569.01'
$ws.Range('E29').Value = '  -3.43%  '
Set-TextValue $ws.Range('D30') '519.67'
$ws.Range('E30').Value = '  -3.34%  '
Set-TextValue $ws.Range('D31') '0.0₃0908'
$ws.Range('E31').Value = '  -6.31%  '
$ws.Range('E32').Value = '  -8.18%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D33') '1.24'
$ws.Range('E33').Value = '  -6.00%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D34') '1.78'
$ws.Range('E34').Value = '  -5.60%  '
$ws.Range('E35').Value = '  +0.08%  '
Set-TextValue $ws.Range('D36') '0.119'
$ws.Range('E36').Value = '  -7.58%  '
Set-TextValue $ws.Range('D37') '157.78'
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  -3.83%  '
$ws.Range('E40').Value = '  -6.08%  '
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D42') '1.67'
$ws.Range('E42').Value = '  -6.45%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range('D43') '4.81'
$ws.Range('E43').Value = '  -6.72%  '
$ws.Range('E44').Value = '  -7.04%  '
Set-TextValue $ws.Range('D45') '2.39'
$ws.Range('E45').Value = '  -7.37%  '
Set-TextValue $ws.Range('D46') '38.71'
$ws.Range('E46').Value = '  -0.94%  '
Set-TextValue $ws.Range('D47') '141.07'
$ws.Range('E47').Value = '  -4.08%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D48') '0.516'
$ws.Range('E48').Value = '  -6.99%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D49') '3.46'
$ws.Range('E49').Value = '  -6.94%  '
$ws.Range('E50').Value = '  -12.19%  '
Set-TextValue $ws.Range('D51') '1.60'
$ws.Range('E51').Value = '  -7.35%  '
